$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8 and 9 contain duplicate "boy/zipper" sentences that were left behind
# when the "man/brochure" and "man/gate" sentence groups were inserted above.
# Delete them entirely so the following rows (child/ketchup bottle, etc.) shift
# up to fill the gap.
$ws.Rows("8:9").Delete()

$ws.Range("A2").Select()
